# Update "想去人数" (want-to-go count) values in both the "展览" and
# "全部类型" sheets to reflect newly scraped figures.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 10468
$ws1.Range("F18").Value = 181
$ws1.Range("F20").Value = 3286

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 10468
$ws4.Range("F21").Value = 181
$ws4.Range("F23").Value = 3286
